# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F, shared by both sheets.
$updates = @{
    2  = 68
    3  = 1050
    4  = 44
    6  = 2968
    8  = 1996
    9  = 187
    10 = 102
    11 = 835
    12 = 32
    13 = 29
    14 = 213
    17 = 31
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
